$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for the rows whose data block
# (columns D, K, L, M, N, O, P, Q, R, S, T) gets cyclically shifted.
$rows = @(2, 3, 4, 6, 7, 8)

$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# New value for row[i] comes from the old value of the next row in the
# cycle (2 -> 3 -> 4 -> 6 -> 7 -> 8 -> 2).
$count = $rows.Count
for ($i = 0; $i -lt $count; $i++) {
    $target = $rows[$i]
    $source = $rows[($i + 1) % $count]
    $src = $data[$source]

    $ws.Cells.Item($target, 4).Value2 = $src.D
    $ws.Cells.Item($target, 11).Value2 = $src.K
    $ws.Cells.Item($target, 12).Value2 = $src.L
    $ws.Cells.Item($target, 13).Value2 = $src.M
    $ws.Cells.Item($target, 14).Value2 = $src.N
    $ws.Cells.Item($target, 15).Value2 = $src.O
    $ws.Cells.Item($target, 16).Value2 = $src.P
    $ws.Cells.Item($target, 17).Value2 = $src.Q
    $ws.Cells.Item($target, 18).Value2 = $src.R
    $ws.Cells.Item($target, 19).Value2 = $src.S
    $ws.Cells.Item($target, 20).Value2 = $src.T
}
